$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: Team Hours / Cumulative Hours / Week Earned Value / Cumulative Earned Value
$ws.Range("F18").Value = 100
$ws.Range("G18").Value = 100
$ws.Range("H18").Value = $ws.Range("E18").Value2
$ws.Range("I18").Value = $ws.Range("E18").Value2

# Row 19
$ws.Range("F19").Value = 120
$ws.Range("G19").Value = 220
$ws.Range("H19").Value = 0.2609
$ws.Range("I19").Value = $ws.Range("E19").Value2

# Apply percentage number format (matches style index 9 / numFmtId 10 "0.00%") to H18:I19
$ws.Range("H18:I19").NumberFormat = "0.00%"

# Update the active selection to I23
$ws.Range("I23").Select()
